# Reproduces the authoring commit that:
#   1. Re-applies a built-in table style to the single table in the deck
#      (slide 16): {0E8FFC6B-185B-4FD9-8F19-A91F3D9FB7B2} (default "no
#      style, no grid") -> {111AF250-6720-4427-AA3D-5C0294D188D1}.
#   2. Re-applies the slide master / notes master design themes ("Integral"
#      and "Office Theme" respectively) -- this is a purely cosmetic,
#      no-observable-effect re-application (same relationships afterwards),
#      included here for completeness/fidelity with the original action.

$p = $ppt.ActivePresentation

# --- 1. Locate the only table shape in the presentation and re-apply its
#        style via the Table Design gallery equivalent (Table.ApplyStyle).
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle("{111AF250-6720-4427-AA3D-5C0294D188D1}", $true)
        }
    }
}

# --- 2. Re-apply (refresh) the design themes already in use by the slide
#        master and the notes master.
$p.SlideMaster.ApplyTheme($p.SlideMaster.Name)
$p.NotesMaster.ApplyTheme($p.NotesMaster.Name)
